$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are removed entirely in the target version
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update values that changed due to forecaster recalculation
$ws.Range("E3").Value = 19.1981274365808
$ws.Range("E4").Value = 0.5079568386449518
$ws.Range("E5").Value = -9.964084247724713
$ws.Range("C6").Value = -0.9140166223623458
$ws.Range("E6").Value = 5.639535270494123
$ws.Range("C7").Value = -5.440152375872276
$ws.Range("E7").Value = -14.43639438706736
$ws.Range("C9").Value = 9.349082908138474
$ws.Range("E9").Value = 27.15801420548431
$ws.Range("E10").Value = -0.6955733540840225
$ws.Range("C11").Value = 0.5389546843749926
$ws.Range("C13").Value = -4.232836797447703
$ws.Range("E13").Value = -8.093075920532211
$ws.Range("C14").Value = -1.479696720105184
$ws.Range("E15").Value = 15.37760125310901
$ws.Range("E16").Value = -2.092856741436244
$ws.Range("C17").Value = 1.913895196850146
$ws.Range("E17").Value = 6.974907992137958
$ws.Range("E18").Value = -1.259568900987029
$ws.Range("C19").Value = 4.861901970954019
$ws.Range("C20").Value = 2.192778679161966
$ws.Range("C21").Value = 4.115488239647735
$ws.Range("E21").Value = 9.52305004616103
$ws.Range("C22").Value = 3.408364488606752
$ws.Range("E22").Value = -0.3010260522302244
$ws.Range("C23").Value = -0.1170030518464982
$ws.Range("E23").Value = -1.421675245489573
$ws.Range("C24").Value = 4.073811422566442
$ws.Range("E24").Value = 8.549566886636839
$ws.Range("E26").Value = 3.825329033908798
$ws.Range("C27").Value = 4.462957198457018
$ws.Range("E28").Value = 7.819356632099983
$ws.Range("C29").Value = 3.65644131722509
$ws.Range("E29").Value = -3.202316982060582
$ws.Range("E31").Value = 2.423519345863312
$ws.Range("E33").Value = -13.80331328666086
$ws.Range("C34").Value = 1.666553973046025
$ws.Range("E34").Value = -7.134843267358049
$ws.Range("E35").Value = 7.549216998964559
$ws.Range("C36").Value = 3.285232806602423
$ws.Range("E36").Value = 4.648946574958668
$ws.Range("C37").Value = 1.566171461167398
$ws.Range("E37").Value = 0.9875809772451616
$ws.Range("C38").Value = 1.879266440112781
$ws.Range("E38").Value = -6.539839435602913
$ws.Range("C39").Value = -1.781446606796444
$ws.Range("E41").Value = -11.32548012975855
$ws.Range("C42").Value = -2.620683231370935
$ws.Range("E42").Value = -5.419975784955122
$ws.Range("C44").Value = 1.16693824877212
$ws.Range("E44").Value = 16.68718678695833
$ws.Range("C45").Value = -2.536896655888543
$ws.Range("C46").Value = -3.036556262700263
$ws.Range("C47").Value = -4.691918671600925
$ws.Range("E47").Value = -6.71087780539289
$ws.Range("C48").Value = 1.758584501904181
$ws.Range("C49").Value = -3.187937113772665
$ws.Range("E49").Value = -7.130744761216423
$ws.Range("E51").Value = 4.240665037474822
$ws.Range("C52").Value = 0.5908161348962437
$ws.Range("E53").Value = -7.500482081224535
